# Update the public EPEX Spot prices workbook:
# Add a new column "I" with the prices for 22-jun on the "Prix Spot" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Header cell (I1): same text + formatting as the other header cells (H1).
$ws.Range("I1").Value = "22-jun"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# Data values for column I (rows 2-25).
$values = @(
    77.02,
    59.91,
    39.84,
    30.33,
    25,
    27.95,
    28.27,
    26.44,
    0,
    -0.09,
    -11.11,
    -33.1,
    -51.14,
    -78.05,
    -64.34999999999999,
    -32.8,
    -7.49,
    0,
    27.56,
    71.14,
    114.12,
    115.32,
    112.43,
    101.63
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 9).Value = $v
    $row = $row + 1
}
